$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.943.59'
$ws.Range("E2").Value = '  -0.56%  '
$ws.Range("D3").Value = '2.218.94'
$ws.Range("E3").Value = '  -0.79%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '292.15'
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").Value = '86.82'
$ws.Range("E6").Value = '  -2.08%  '
$ws.Range("E7").Value = '  -0.92%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.466'
$ws.Range("D10").Value = '30.45'
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("D11").Value = '50.39'
$ws.Range("E11").Value = '  +5.61%  '
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("E13").Value = '  +3.16%  '
$ws.Range("D14").Value = '6.43'
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("D15").Value = '2.561.88'
$ws.Range("E15").Value = '  -0.62%  '
$ws.Range("D16").Value = '13.79'
$ws.Range("E16").Value = '  -3.05%  '
$ws.Range("D17").Value = '2.206.90'
$ws.Range("E17").Value = '  -2.32%  '
$ws.Range("D18").Value = '0.732'
$ws.Range("E18").Value = '  -0.80%  '
$ws.Range("D19").Value = '39.874.62'
$ws.Range("E19").Value = '  -0.61%  '
$ws.Range("E20").Value = '  -0.67%  '
$ws.Range("E21").Value = '  -4.74%  '
$ws.Range("D22").Value = '5.75'
$ws.Range("E22").Value = '  -2.12%  '
$ws.Range("D23").Value = '65.54'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").Value = '237.04'
$ws.Range("E24").Value = '  +0.16%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("D26").Value = '2.47'
$ws.Range("E26").Value = '  -0.68%  '
$ws.Range("D27").Value = '1.82'
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("E28").Value = '  +7.57%  '
$ws.Range("D29").Value = '23.29'
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = '9.23'
$ws.Range("E30").Value = '  -1.29%  '
$ws.Range("D31").Value = '157.50'
$ws.Range("E31").Value = '  +3.19%  '
$ws.Range("D32").Value = '31.80'
$ws.Range("E32").Value = '  -4.09%  '
$ws.Range("E33").Value = '  +0.08%  '
$ws.Range("D34").Value = '4.97'
$ws.Range("E34").Value = '  -0.62%  '
$ws.Range("E35").Value = '  +3.72%  '
$ws.Range("D36").Value = '0.0714'
$ws.Range("E36").Value = '  -1.34%  '
$ws.Range("E37").Value = '  -1.37%  '
$ws.Range("E38").Value = '  -0.71%  '
$ws.Range("D39").Value = '0.0989'
$ws.Range("E39").Value = '  -1.50%  '
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("D41").Value = '15.18'
$ws.Range("E41").Value = '  -6.65%  '
$ws.Range("D42").Value = '2.085.10'
$ws.Range("E42").Value = '  -1.26%  '
$ws.Range("E43").Value = '  -3.73%  '
$ws.Range("E44").Value = '  -0.49%  '
$ws.Range("D45").Value = '17.94'
$ws.Range("E45").Value = '  -1.58%  '
$ws.Range("D46").Value = '9.79'
$ws.Range("E46").Value = '  -3.05%  '
$ws.Range("D47").Value = '1.98'
$ws.Range("E47").Value = '  -9.07%  '
$ws.Range("E48").Value = '  +0.21%  '
$ws.Range("D49").Value = '2.434.46'
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("E51").Value = '  +1.69%  '
